$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("B9").Value = 197500
$ws.Range("C9").Value = 1
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = ";0"
$ws.Range("G9").Value = ";12"
$ws.Range("H9").Value = ";-500"
